$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "KNN_version_3"
$ws.Range("B4").Value = 0.3217391304347826
$ws.Range("C4").Value = 0.07035330694830499
$ws.Range("D4").Value = 0.07430097822851121
$ws.Range("E4").Value = 0.07126649005254523
$ws.Range("F4").Value = 0.2836746693384991
$ws.Range("G4").Value = 0.3217391304347826
$ws.Range("H4").Value = 0.2985569894858061

$ws.Range("A5").Value = "KNN_version_4"
$ws.Range("B5").Value = 0.2363636363636364
$ws.Range("C5").Value = 0.1128668149459225
$ws.Range("D5").Value = 0.1150636802810716
$ws.Range("E5").Value = 0.113201041836724
$ws.Range("F5").Value = 0.2304672457816474
$ws.Range("G5").Value = 0.2363636363636364
$ws.Range("H5").Value = 0.2323304717634825
